# Apply BoM/Costs updates: R4 and U2 removed, leaving only R2 and U1.
# This reduces component/quantity counts across the BoM, DNF and Costs sheets,
# and refreshes the "Created" timestamp.

$wb = $excel.ActiveWorkbook

$wsBoM      = $wb.Worksheets.Item("BoM")
$wsDNF      = $wb.Worksheets.Item("DNF")
$wsCosts    = $wb.Worksheets.Item("Costs")
$wsCostsDNF = $wb.Worksheets.Item("Costs (DNF)")

# --- BoM sheet -----------------------------------------------------------
# References column: "R2 R4" -> "R2", "U1 U2" -> "U1"
$wsBoM.Range("D16").Value = "R2"
$wsBoM.Range("D18").Value = "U1"

# Quantity Per PCB column: 2 -> 1 for both rows
$wsBoM.Range("G16").Value = 1
$wsBoM.Range("G18").Value = 1

# Header summary cells
$wsBoM.Range("F3").Value = "216 (205 SMD/ 0 THT)"
$wsBoM.Range("F4").Value = "213 (204 SMD/ 0 THT)"

# Total Components summary cell
$wsBoM.Range("F6").Value = 213

# --- DNF sheet -------------------------------------------------------------
# Header summary cells
$wsDNF.Range("F3").Value = "216 (205 SMD/ 0 THT)"
$wsDNF.Range("F4").Value = "213 (204 SMD/ 0 THT)"

# Total Components summary cell
$wsDNF.Range("F6").Value = 213

# --- Costs sheet -----------------------------------------------------------
# References column: "R2 R4" -> "R2", "U1 U2" -> "U1"
$wsCosts.Range("A17").Value = "R2"
$wsCosts.Range("A19").Value = "U1"

# Build Quantity formulas: CEILING(BoardQty*2,1) -> BoardQty*1
$wsCosts.Range("F17").Formula = "=BoardQty*1"
$wsCosts.Range("F19").Formula = "=BoardQty*1"

# Created timestamp
$wsCosts.Range("B22").Value = "2024-10-05 16:27:13"

# --- Costs (DNF) sheet ------------------------------------------------------
# Created timestamp
$wsCostsDNF.Range("B14").Value = "2024-10-05 16:27:13"
